# Updated cryptos list on Tue Jan 16 09:40:09 UTC 2024 with GitHub Actions
#
# Price (column D) and Volume/1h (column E) figures are refreshed for each
# coin row. Column D holds plain text (not numbers) in this workbook, even
# for values that look numeric - e.g. thousand-grouped prices like
# "43.083.61" are never valid floats, but plain decimals like "318.82" are,
# and Excel would silently coerce those to a numeric cell on assignment.
# To keep every D-cell a verbatim text value (matching the original
# t="inlineStr" cells, with no numeric rounding/precision drift and no
# left-over number-format style), values that parse as a float are written
# with a leading apostrophe (forces text entry) and the cell style is then
# reset back to "Normal" to drop the resulting quote-prefix formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.083.61"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.545.37"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'318.82"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'98.06"
$ws.Range("E6").Value = "  +3.12%  "
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'36.42"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").Value = "'0.0819"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").Value = "2.934.45"
$ws.Range("D15").Value = "2.573.04"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "'15.23"
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "'0.854"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "43.094.49"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "'12.91"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "0.0₃0971"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'70.09"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "'255.16"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'2.97"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'2.07"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("D26").Value = "'26.69"
$ws.Range("E26").Value = "  -3.53%  "
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("D29").Value = "'41.05"
$ws.Range("E29").Value = "  +4.35%  "
$ws.Range("E30").Value = "  +4.19%  "
$ws.Range("D31").Value = "'5.96"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("D32").Value = "'158.15"
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("D33").Value = "'2.19"
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").Value = "  +4.28%  "
$ws.Range("D36").Value = "'19.10"
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "'2.48"
$ws.Range("E39").Value = "  +13.14%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").Value = "'22.18"
$ws.Range("E41").Value = "  -9.35%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'0.0305"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "'3.31"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "2.022.66"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").Value = "'9.15"
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("D48").Value = "'84.75"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").Value = "'76.77"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("D50").Value = "'106.85"
$ws.Range("E50").Value = "  +4.70%  "
$ws.Range("D51").Value = "2.791.00"
$ws.Range("E51").Value = "  +0.83%  "

# Strip the quote-prefix formatting picked up above so these cells end up
# style-identical to their neighbours (no "s" attribute), same as source.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
